$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row 3 value - plain text (shared string), not a hyperlink
$ws.Range("A3").Value = "https://ingenarte.github.io/react-tetris2/"

# Update selection to A7 as recorded in the saved view state
$ws.Range("A7").Select()
